# follow_map_position.xlsx - translations check-in
#
# Content changes:
#   - survey!F1   "display.text"   -> "display.prompt.text"
#   - settings!C1 "display.title"  -> "display.title.text"
#
# View/selection changes:
#   - survey      becomes the active (tabSelected) sheet, selection -> F2
#   - settings    selection -> C2 (no longer special)
#   - properties  loses tabSelected (was active before), selection stays E7

$wb = $excel.ActiveWorkbook

$survey     = $wb.Worksheets.Item("survey")
$settings   = $wb.Worksheets.Item("settings")
$properties = $wb.Worksheets.Item("properties")

# -- text updates (order matters for shared-string table ordering) --
$settings.Range("C1").Value = "display.title.text"
$survey.Range("F1").Value = "display.prompt.text"

# -- re-establish selections on the non-active sheets first --
$settings.Range("C2").Select()
$properties.Range("E7").Select()

# -- finally select survey!F2, leaving it the active sheet/tab --
$survey.Range("F2").Select()
